$d = $word.ActiveDocument

# Title paragraph: "Questions:" | " " | "Arithmetic" | " " | "on" | " " |
# "complex" | " " | "numbers"  ->  single run "Questions: Arithmetic on
# complex numbers". Scope the Find/Replace to the paragraph's own Range
# so we don't touch identical-looking text elsewhere in the document.
$pTitle = $d.Paragraphs(1)
$pTitle.Range.Find.Execute(
    "Questions: Arithmetic on complex numbers", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Questions: Arithmetic on complex numbers", 2)

# Author paragraph: "Charlotte" | " " | "McCarthy" -> "Charlotte McCarthy"
$pAuthor = $d.Paragraphs(2)
$pAuthor.Range.Find.Execute(
    "Charlotte McCarthy", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Charlotte McCarthy", 2)

# Abstract paragraph: many single-word/space runs -> one run with the
# full sentence.
$pAbstract = $d.Paragraphs(4)
$pAbstract.Range.Find.Execute(
    "A selection of questions for the study guide on arithmetic on complex numbers.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A selection of questions for the study guide on arithmetic on complex numbers.",
    2)
